$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the missing "Hours" value for the existing 3rd data row (row 4)
$ws.Range("B4").Value = 2

# Add a brand-new timesheet row (row 5)
$ws.Range("B5").Value = 3.25
$ws.Range("C5").Value = "6:45pm"
$ws.Range("D5").Value = "10pm"

# Move the active selection the way Excel does after data entry
$ws.Range("C7").Select()
